$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the twentieth lesson (row 22)
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 43136
$ws.Cells.Item(22, 2).Style = $ws.Cells.Item(21, 2).Style
$ws.Cells.Item(22, 3).Value = 2

# Add the twenty-first lesson (row 23)
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 43138
$ws.Cells.Item(23, 2).Style = $ws.Cells.Item(21, 2).Style
$ws.Cells.Item(23, 3).Value = 2

# Align D14's style with D10 (italic "New Set of Hours" note)
$ws.Range("D14").Style = $ws.Range("D10").Style

# Move the active selection
$ws.Range("A24").Select()
